# Add files via upload
#
# Fills in Outcome / Cause of Action / Civil or Criminal / Token / Project Name /
# Blockchain / Amount / Securities Act of 1933 / Securities Exchange Act of 1934 /
# SEC Office data for rows 49-54, which previously only had Date/Case/Description
# (columns A-C) populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 49 : SEC v. Grybniak, et al. (Opporty International / OPP) ----
$ws.Range("D49").Value = "Ongoing"
$ws.Range("E49").Value = "Unregistered Offering"
$ws.Range("F49").Value = "Civil"
$ws.Range("I49").Value = "Ethereum"
$ws.Range("J49").Value = 600000
$ws.Range("K49").Value = 1
$ws.Range("L49").Value = 1
$ws.Range("M49").Value = "Washington, D.C."

# ---- Row 50 : Blockchain of Things, Inc. (BCOT) ----
$ws.Range("D50").Value = "Settlement"
$ws.Range("E50").Value = "Unregistered Offering"
$ws.Range("F50").Value = "Civil"
$ws.Range("I50").Value = "N/A"
$ws.Range("J50").Value = 13000000
$ws.Range("K50").Value = 1
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = "Washington, D.C."

# ---- Row 51 : SEC v. Eyal, et al. (UnitedData, Inc. / Shopin) ----
$ws.Range("D51").Value = "Ongoing"
$ws.Range("E51").Value = "Unregistered Offering"
$ws.Range("F51").Value = "Civil"
$ws.Range("I51").Value = "Ethereum"
$ws.Range("J51").Value = 42500000
$ws.Range("K51").Value = 1
$ws.Range("L51").Value = 1
$ws.Range("M51").Value = "New York"

# ---- Row 52 : SEC v. Telegram Group Inc., et al. (Grams / TON) ----
$ws.Range("D52").Value = "Settlement"
$ws.Range("E52").Value = "Unregistered Offering"
$ws.Range("F52").Value = "Civil"
$ws.Range("J52").Value = 1700000000
$ws.Range("K52").Value = 1
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = "New York"

# ---- Row 53 : Block.one (EOS) ----
$ws.Range("D53").Value = "Settlement"
$ws.Range("E53").Value = "Unregistered Offering"
$ws.Range("F53").Value = "Civil"
$ws.Range("H53").Value = "Block.one"
$ws.Range("I53").Value = "Ethereum"
$ws.Range("J53").Value = 24000000
$ws.Range("K53").Value = 1
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = "New York"

# ---- Row 54 : SEC v. Lucas (Fantasy Market / FMT) ----
$ws.Range("D54").Value = "Settlement"
$ws.Range("E54").Value = "Unregistered Offering and Fraud"
$ws.Range("F54").Value = "Civil"
$ws.Range("I54").Value = "Ethereum"
$ws.Range("J54").Value = 63000
$ws.Range("K54").Value = 1
$ws.Range("L54").Value = 1
$ws.Range("M54").Value = "New York"

# ---- Token / Project Name cells (new vocabulary) ----
$ws.Range("H49").Value = "Opporty International, Inc."
$ws.Range("G49").Value = "OPP"
$ws.Range("G50").Value = "BCOT"
$ws.Range("H51").Value = "UnitedData, Inc."
$ws.Range("H50").Value = "Blockchain of Things Inc. "
$ws.Range("G51").Value = "Shopin"
$ws.Range("H52").Value = "Telegram Group Inc."
$ws.Range("G52").Value = "Grams"
$ws.Range("I52").Value = "TON"
$ws.Range("G53").Value = "EOS"
$ws.Range("G54").Value = "FMT"
$ws.Range("H54").Value = "Fantasy Market"

# ---- Restore the author's final scroll/selection state ----
$ws.Range("F55").Select()
